$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Prime the formatting of the new rows by copying from existing, same-banding rows ---
# Row 835 continues the single-row style band (fill "10"/"44", like row 821).
$ws.Range("A821:E821").Copy()
$ws.Range("A835:E835").PasteSpecial(-4122)

# Rows 836:843 start a new style band (fill "11"/"45", like rows 813:818).
$ws.Range("A813:E818").Copy()
$ws.Range("A836:E841").PasteSpecial(-4122)
$ws.Range("A813:E814").Copy()
$ws.Range("A842:E843").PasteSpecial(-4122)

$ws.Range("A835:E843").RowHeight = 15.75

# --- Row 835: ABDALLH OMAR ABDALLH MUHAMMAD ---
$ws.Range("A835").Value = "DSS1834"
$ws.Range("B835").Value = "ABDALLH OMAR ABDALLH MUHAMMAD"
$ws.Range("C835").Value = "Health and Safety Advanced Diploma"
$ws.Range("D835").Value = "'15-12-2024"
$ws.Range("E835").Value = 1

# --- Rows 836:843: Mahmoud Elsayed Abdellah Mohamed ---
$certNos = @("DSS1835","DSS1836","DSS1837","DSS1838","DSS1839","DSS1840","DSS1841","DSS1842")
$courses = @(
  "30 Hours Construction Safety & Health",
  "30 Hours G. Industry Safety & Health",
  "Electrical Safety & LOTO",
  "Fire Marshal",
  "Scaffold Competent Person",
  "Lifting & Rigging Competent Person",
  "Health & Safety Risk Assessment",
  "Safety Management System & PTW"
)
$dates = @("'05-12-2024","'10-12-2024","'06-12-2024","'03-12-2024","'01-12-2024","'02-12-2024","'07-12-2024","'08-12-2024")

for ($i = 0; $i -lt 8; $i++) {
  $r = 836 + $i
  $ws.Range("A$r").Value = $certNos[$i]
  $ws.Range("B$r").Value = "Mahmoud Elsayed Abdellah Mohamed"
  $ws.Range("C$r").Value = $courses[$i]
  $ws.Range("D$r").Value = $dates[$i]
  $ws.Range("E$r").Value = 1
}

# --- Mark a new manual page break after the last data row (matches the other ~85-row bands) ---
$ws.Rows.Item(844).PageBreak = 1

# --- Restore the view: scrolled down to the new bottom of the table, selection on C849 ---
$excel.ActiveWindow.ScrollRow = 832
$ws.Range("C849").Select()
